$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 273 - this shifts the existing rows 273:292
# down to 274:293 (and bumps the sheet dimension to A1:R293).
$ws.Rows(273).Insert()

# Populate the newly inserted row 273 with the new weekly record.
# (Columns A,B,C,E,F,G,H,N,O,Q,R repeat the same constant values used
# throughout this Cilantro / Terminal Hortofrutícola Agro Chillán block.)
$ws.Range("A273").Value = 7
$ws.Range("B273").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C273").Value = "Ñuble"
$ws.Range("D273").Value = 45132
$ws.Range("E273").Value = 16
$ws.Range("F273").Value = 100112040
$ws.Range("G273").Value = "Cilantro"
$ws.Range("H273").Value = "Sin especificar"
$ws.Range("I273").Value = "Primera"
$ws.Range("J273").Value = 200
$ws.Range("K273").Value = 1500
$ws.Range("L273").Value = 1500
$ws.Range("M273").Value = 1500
$ws.Range("N273").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O273").Value = "Provincia de Diguillín"
$ws.Range("P273").Value = 1500
$ws.Range("Q273").Value = 1
$ws.Range("R273").Value = "Hortaliza"
